$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts old D->E, old E->F), making room for
# the new "subtype" column.
$ws.Columns("D").Insert()

# Match column C's width on the newly inserted subtype column (closest value
# the COM width quantization can reach to the authored 15.1640625).
$ws.Range("D1").ColumnWidth = 14.28

# --- Header row ---
$ws.Range("D1").Value = "subtype"

# --- Row 2 (type markers) ---
$ws.Range("D2").Value = "string"

# --- Row 3 (noSignUp) ---
$ws.Range("C3").Value = "city"
$ws.Range("D3").Value = "percentage"

# --- Row 4 (haveRecommendLetter) ---
$ws.Range("D4").Value = "item"

# --- Row 5 (canSpreadRumor) ---
$ws.Range("D5").Value = "job"

# --- guild type (rows 4 & 5) ---
$ws.Range("C4").Value = "guild"
$ws.Range("C5").Value = "guild"

# --- compareType / parameter updates ---
$ws.Range("E3").Value = "'<="
$ws.Range("F3").Value = 20
$ws.Range("F5").Value = 12

# Update the selection to match the authored view state.
$ws.Range("F3").Select()
